{"js": "// Replace each two-digit-by-two-digit multiplication prompt text with its\n// new pair of operands. The mapping below is applied strictly in document\n// order (top-left to bottom-right of the table) so that a new value that\n// happens to equal another cell's *original* value (e.g. \"65\u00d738=\" is both\n// an original prompt earlier in the table and the replacement text used\n// later) never gets re-matched by a later search.\nconst replacements = [\n  [\"30\u00d711=\", \"89\u00d797=\"],\n  [\"22\u00d727=\", \"49\u00d777=\"],\n  [\"95\u00d737=\", \"50\u00d751=\"],\n  [\"65\u00d727=\", \"87\u00d758=\"],\n  [\"31\u00d797=\", \"28\u00d789=\"],\n  [\"77\u00d738=\", \"66\u00d786=\"],\n  [\"63\u00d782=\", \"21\u00d771=\"],\n  [\"19\u00d782=\", \"46\u00d711=\"],\n  [\"69\u00d757=\", \"72\u00d748=\"],\n  [\"57\u00d770=\", \"63\u00d725=\"],\n  [\"17\u00d749=\", \"84\u00d726=\"],\n  [\"71\u00d750=\", \"78\u00d747=\"],\n  [\"19\u00d721=\", \"18\u00d740=\"],\n  [\"65\u00d738=\", \"69\u00d769=\"],\n  [\"99\u00d773=\", \"71\u00d784=\"],\n  [\"92\u00d723=\", \"86\u00d712=\"],\n  [\"15\u00d752=\", \"66\u00d746=\"],\n  [\"28\u00d739=\", \"71\u00d714=\"],\n  [\"90\u00d764=\", \"34\u00d731=\"],\n  [\"50\u00d775=\", \"51\u00d712=\"],\n  [\"43\u00d733=\", \"43\u00d717=\"],\n  [\"42\u00d784=\", \"65\u00d738=\"],\n  [\"70\u00d728=\", \"47\u00d761=\"],\n  [\"32\u00d727=\", \"59\u00d759=\"],\n  [\"31\u00d784=\", \"44\u00d736=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  // Only the first match is relevant: each search string is unique in the\n  // document at the moment it is searched for (see ordering note above).\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-by-two-digit multiplication prompt text with its\n# new pair of operands. The mapping below is applied strictly in document\n# order (top-left to bottom-right of the table) so that a new value that\n# happens to equal another cell's *original* value (e.g. \"65\u00d738=\" is both\n# an original prompt earlier in the table and the replacement text used\n# later) never gets re-matched by a later search: each \"find\" string is\n# unique in the document at the moment it is searched for.\n$pairs = @(\n    @(\"30\u00d711=\", \"89\u00d797=\"),\n    @(\"22\u00d727=\", \"49\u00d777=\"),\n    @(\"95\u00d737=\", \"50\u00d751=\"),\n    @(\"65\u00d727=\", \"87\u00d758=\"),\n    @(\"31\u00d797=\", \"28\u00d789=\"),\n    @(\"77\u00d738=\", \"66\u00d786=\"),\n    @(\"63\u00d782=\", \"21\u00d771=\"),\n    @(\"19\u00d782=\", \"46\u00d711=\"),\n    @(\"69\u00d757=\", \"72\u00d748=\"),\n    @(\"57\u00d770=\", \"63\u00d725=\"),\n    @(\"17\u00d749=\", \"84\u00d726=\"),\n    @(\"71\u00d750=\", \"78\u00d747=\"),\n    @(\"19\u00d721=\", \"18\u00d740=\"),\n    @(\"65\u00d738=\", \"69\u00d769=\"),\n    @(\"99\u00d773=\", \"71\u00d784=\"),\n    @(\"92\u00d723=\", \"86\u00d712=\"),\n    @(\"15\u00d752=\", \"66\u00d746=\"),\n    @(\"28\u00d739=\", \"71\u00d714=\"),\n    @(\"90\u00d764=\", \"34\u00d731=\"),\n    @(\"50\u00d775=\", \"51\u00d712=\"),\n    @(\"43\u00d733=\", \"43\u00d717=\"),\n    @(\"42\u00d784=\", \"65\u00d738=\"),\n    @(\"70\u00d728=\", \"47\u00d761=\"),\n    @(\"32\u00d727=\", \"59\u00d759=\"),\n    @(\"31\u00d784=\", \"44\u00d736=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
